$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rule R30: update the "Integer min" value (cell C10) from 18 to 1
$ws.Range("C10").Value = 1
